# Diary.xlsx — add LeetCode problems 19, 25, 29 & 31 (commit: "Problem 19, 22, 24, 29 & 31")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: Remove Nth Node From End of List (Problem 19) ---------------
# Reuse row 51's formatting (C -> Neutral/no value, D -> Good) for the new row.
$ws.Range("A51:D51").Copy()
$ws.Range("A53:D53").PasteSpecial(-4122)
$ws.Cells.Item(53, 1).Value = 19
$ws.Cells.Item(53, 2).Value = "Remove Nth Node From End of List"

# --- Row 54: Swap Nodes in Pairs (Problem 25) -----------------------------
# Reuse row 52's formatting (C -> Neutral/no value, D -> Neutral).
$ws.Range("A52:D52").Copy()
$ws.Range("A54:D54").PasteSpecial(-4122)
$ws.Cells.Item(54, 1).Value = 25
$ws.Cells.Item(54, 2).Value = "Swap Nodes in Pairs"

# --- Row 55: Divide Two Integers (Problem 29) -----------------------------
# Reuse row 45's formatting (C -> Neutral/no value, D -> Bad).
$ws.Range("A45:D45").Copy()
$ws.Range("A55:D55").PasteSpecial(-4122)
$ws.Cells.Item(55, 1).Value = 29
$ws.Cells.Item(55, 2).Value = "Divide Two Integers"

# --- Row 56: Next Permutation (Problem 31) --------------------------------
# Reuse row 45's formatting (C -> Neutral/no value, D -> Bad).
$ws.Range("A45:D45").Copy()
$ws.Range("A56:D56").PasteSpecial(-4122)
$ws.Cells.Item(56, 1).Value = 31
$ws.Cells.Item(56, 2).Value = "Next Permutation"

$excel.CutCopyMode = 0

# --- Update the view to match where the workbook was scrolled/selected ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D56").Select()
